$d = $word.ActiveDocument

function Insert-AnswerXml {
    param(
        [int]$ParaIndex,
        [string]$InnerXml
    )
    $p = $d.Paragraphs.Item($ParaIndex)
    $r = $p.Range
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $InnerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData>' +
        '</pkg:part>' +
        '</pkg:package>'
    $r.InsertXML($pkg)
}

# --- Bottom of the document first, so earlier paragraph indices stay valid ---

# Paragraph after "The story is about the workers... leadership in the packing house?"
# Previously held the lone _GoBack bookmark; replace with the answer text and drop the bookmark.
Insert-AnswerXml 30 '<w:p><w:r><w:t xml:space="preserve">It’s extremely transactional. The business owners do not care about the workers because there are thousands more lined up outside. </w:t></w:r></w:p>'

# Paragraph after "Look for examples of Sinclair’s use of mechanical language..."
# becomes two paragraphs, each with one quotation.
Insert-AnswerXml 28 '<w:p><w:r><w:t>…they had been serving as cogs in the great packing machine, and now was the time for the renovating of it, and the replacing of damaged parts.</w:t></w:r></w:p><w:p><w:r><w:t>…a little nearer to the time when it would be their turn to be shaken from the tree.</w:t></w:r></w:p>'

# Paragraph after "How does Lao Tzu’s Tao set the groundwork for the concept of Servant Leadership?"
# gets the answer text plus the relocated _GoBack bookmark.
Insert-AnswerXml 17 '<w:p><w:r><w:t xml:space="preserve">Lao thinks that the leader should put the group’s priorities above his own. Leaders should look to remove blocks so that the group can accomplish their goals. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

# Paragraph after "To Lao Tzu, is the leader’s role with the group a passive or active one?"
Insert-AnswerXml 15 '<w:p><w:r><w:t>Passive. See #58 “The less a leader does and says, the happier his people.”</w:t></w:r></w:p>'

# Paragraph after "What was the leader’s relationship with the group according to Lao Tzu?"
Insert-AnswerXml 13 '<w:p><w:r><w:t xml:space="preserve">The leader should honor the group, not the other way around. Leaders should not be looking for admiration. </w:t></w:r></w:p>'

# Paragraph after "In your own words, how did Lao Tzu view the role of the leader?"
Insert-AnswerXml 11 '<w:p><w:r><w:t xml:space="preserve">Similar to a servant-leader. Lao thinks leaders are at best when nobody notices them. </w:t></w:r></w:p>'
